$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2025-10-25 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-26 Sunday", 2)

$t = $d.Tables.Item(1)

# Row 1
$t.Cell(1,1).Range.Text = "154×6=924"
$t.Cell(1,2).Range.Text = "483×9=4347"
$t.Cell(1,3).Range.Text = "537×7=3759"
$t.Cell(1,4).Range.Text = "660×6=3960"
$t.Cell(1,5).Range.Text = "914×4=3656"

# Row 5
$t.Cell(5,1).Range.Text = "843×2=1686"
$t.Cell(5,2).Range.Text = "661×6=3966"
$t.Cell(5,3).Range.Text = "837×4=3348"
$t.Cell(5,4).Range.Text = "230×7=1610"
$t.Cell(5,5).Range.Text = "123×9=1107"

# Row 10
$t.Cell(10,1).Range.Text = "856×6=5136"
$t.Cell(10,2).Range.Text = "137×3=411"
$t.Cell(10,3).Range.Text = "458×9=4122"
$t.Cell(10,4).Range.Text = "615×6=3690"
$t.Cell(10,5).Range.Text = "510×4=2040"

# Row 15
$t.Cell(15,1).Range.Text = "660×5=3300"
$t.Cell(15,2).Range.Text = "642×9=5778"
$t.Cell(15,3).Range.Text = "857×5=4285"
$t.Cell(15,4).Range.Text = "117×8=936"
$t.Cell(15,5).Range.Text = "447×2=894"

# Row 20
$t.Cell(20,1).Range.Text = "665×2=1330"
$t.Cell(20,2).Range.Text = "582×8=4656"
$t.Cell(20,3).Range.Text = "910×6=5460"
$t.Cell(20,4).Range.Text = "416×4=1664"
$t.Cell(20,5).Range.Text = "663×3=1989"
